# Week 13 logging update
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this week's rush/pass yardage entries to the running,
# space-separated per-game logs kept as text in B2:C3.
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value() + " 7 6 0 5 9 2 2 9 4 9 -2 3 0 3 6 8 6 -1 6 7 5 0 25 2 3 2"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value() + " 2 1 3 1 -1 10 0 2 9 -1 7 1 1 -3 13 -1 1 2 -1 4 0 2 1 5 2"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value() + " -1 27 4 35 5 9 8 7 7 7 4 19 8 7 43 29 8 3 19 22 2 5 5 3 8 1"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value() + " 3 28 -1 6 11 5 16 9 -2 12 9 7 17 10 9"

# ---------------------------------------------------------------------------
# OFF sheet: season-to-date offensive totals, updated through Week 13.
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("B2").Value = 4
$wsOFF.Range("C2").Value = 132
$wsOFF.Range("D2").Value = 7
$wsOFF.Range("F2").Value = 31
$wsOFF.Range("G2").Value = 40
$wsOFF.Range("J2").Value = 23
$wsOFF.Range("L2").Value = 232
$wsOFF.Range("M2").Value = 158
$wsOFF.Range("O2").Value = 16
$wsOFF.Range("P2").Value = 8
$wsOFF.Range("Q2").Value = 370

$wsOFF.Range("B3").Value = 6
$wsOFF.Range("C3").Value = 155
$wsOFF.Range("D3").Value = 7
$wsOFF.Range("E3").Value = 17
$wsOFF.Range("F3").Value = 96
$wsOFF.Range("G3").Value = 38
$wsOFF.Range("H3").Value = 9
$wsOFF.Range("I3").Value = 37
$wsOFF.Range("J3").Value = 50
$wsOFF.Range("N3").Value = 9

# ---------------------------------------------------------------------------
# DEF sheet: season-to-date defensive totals, updated through Week 13.
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value = 149
$wsDEF.Range("D2").Value = 8
$wsDEF.Range("F2").Value = 39
$wsDEF.Range("G2").Value = 53
$wsDEF.Range("H2").Value = 7
$wsDEF.Range("I2").Value = 7
$wsDEF.Range("J2").Value = 22
$wsDEF.Range("L2").Value = 219
$wsDEF.Range("M2").Value = 154
$wsDEF.Range("O2").Value = 23
$wsDEF.Range("P2").Value = 11
$wsDEF.Range("Q2").Value = 395

$wsDEF.Range("B3").Value = 11
$wsDEF.Range("C3").Value = 132
$wsDEF.Range("E3").Value = 29
$wsDEF.Range("F3").Value = 88
$wsDEF.Range("G3").Value = 22
$wsDEF.Range("H3").Value = 19
$wsDEF.Range("I3").Value = 52
$wsDEF.Range("J3").Value = 35
$wsDEF.Range("N3").Value = 16

# ---------------------------------------------------------------------------
# ST sheet: special-teams totals plus the per-kicker/returner logs.
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 71
$wsST.Range("D2").Value = 36
$wsST.Range("F2").Value = 100
$wsST.Range("G2").Value = 94
$wsST.Range("J2").Value = 36
$wsST.Range("K2").Value = 33
$wsST.Range("L2").Value = 28
$wsST.Range("M2").Value = 23

$wsST.Range("B3").Value = 42

$wsST.Range("B4").Value = $wsST.Range("B4").Value() + " 58 63"
$wsST.Range("B5").Value = $wsST.Range("B5").Value() + " 24 0"
$wsST.Range("B6").Value = $wsST.Range("B6").Value() + " 65 16"
$wsST.Range("D3").Value = $wsST.Range("D3").Value() + " 46 46 58"
$wsST.Range("D4").Value = $wsST.Range("D4").Value() + " 0 0 0"
$wsST.Range("D5").Value = $wsST.Range("D5").Value() + " 0 0 0 0 0"

# ---------------------------------------------------------------------------
# TURNS sheet: home/road turnover counts.
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("D2").Value = 3
$wsTURNS.Range("E2").Value = 7
$wsTURNS.Range("D3").Value = 2
$wsTURNS.Range("E3").Value = 4

# ---------------------------------------------------------------------------
# PEN sheet: penalty counts.
# ---------------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value = 7
$wsPEN.Range("D3").Value = 5
